$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '51.764.49'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.77%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.996.48'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.30%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '384.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '104.96'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.42%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.51%  '

$ws.Range("E8").Value = '  -0.03%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.600'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '37.52'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.30%  '

$ws.Range("E11").Value = '  +0.59%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0850'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.480.85'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +3.51%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '18.50'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.67%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.63'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +4.20%  '

$ws.Range("E16").Value = '  +11.79%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.005.13'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.41%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '51.812.63'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.83%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.31'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.80%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.02'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.55%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0968'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.20'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '264.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.07%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.95'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +10.15%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.42'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +19.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.71'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +24.35%  '

$ws.Range("E28").Value = '  +2.08%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.115'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +13.79%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '26.13'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.44%  '

$ws.Range("E31").Value = '  +0.06%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.94'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.27%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '35.23'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.94%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.14'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.06%  '

$ws.Range("B35").Value = 'Toncoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.07'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.72%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0454'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.16%  '

$ws.Range("E37").Value = '  -0.25%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.05'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.62%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '17.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.60'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.63%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.85'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.34%  '

$ws.Range("E42").Value = '  +4.46%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '122.90'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '21.94'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.26%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.280'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +19.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.67%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +6.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.042.72'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.65%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0336'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +9.92%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.868'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.09%  '
